$d = $word.ActiveDocument

# --- Locate the "Sandbox" task bullet and the lone "8" run inside it -------
# (".. - 8 SP - Migue"). That "8" is its own run in the OOXML, flanked on
# both sides by sibling runs that share identical run formatting
# (<w:b/>), so a naive text replace would make Word re-coalesce it with
# its neighbours into one big run. We avoid that by bracketing the digit
# with a pair of bookmarks *before* touching the text: a bookmark start
# acts as a hard boundary runs won't be merged across.
$sandboxPara = $null
foreach ($p in $d.Paragraphs) {
  if ($p.Range.Text -like "*Sandbox*") {
    $sandboxPara = $p
    break
  }
}

$digit = $sandboxPara.Range.Duplicate
$digit.Find.Execute("8", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Temporary boundary just to the left of the digit, so the edited run
# won't fuse with the "– " run before it.
$leftEdge = $d.Range($digit.Start, $digit.Start)
$d.Bookmarks.Add("zzTmpBoundary", $leftEdge) | Out-Null

# The real mark Word leaves behind after an edit is the single, special
# "_GoBack" bookmark. (Re)adding it here drops it right after the digit
# and, being unique, automatically supersedes the "_GoBack" bookmark that
# used to sit at the end of the "Publicidad por redes sociales ... Álvaro
# D." paragraph - so that stale one disappears on its own.
$rightEdge = $d.Range($digit.End, $digit.End)
$d.Bookmarks.Add("_GoBack", $rightEdge) | Out-Null

# 8 SP -> 11 SP
$digit.Text = "11"

# Drop the temporary helper bookmark; only "_GoBack" should remain.
$d.Bookmarks("zzTmpBoundary").Delete()
